# Tidsregistrering i PTE projektet Jeppe Damborg.xlsx
# Har rettet Aktivitetsdiagram for UC6 til og opdateret tidsplan for idag
#
# Adds a new day (15-3-2017 / serial 42809) of time-registration entries
# (rows 34-40) plus the trailing total row (41) to the "Tidsregistrering"
# worksheet, mirroring the existing rows' layout/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Copy existing formatting down onto the new rows --------------------
# Column A uses the date style (style index 1 in the original workbook)
$ws.Range("A29").Copy()
$ws.Range("A34").PasteSpecial(-4122)

# Columns G:H use the time style (style index 5 in the original workbook)
$ws.Range("G29:H29").Copy()
$ws.Range("G34:H40").PasteSpecial(-4122)

# --- New day header (row 34) --------------------------------------------
$ws.Range("A34").Value = 42809

# --- Activity rows --------------------------------------------------------
$ws.Range("E34").Value = "Software Architect"
$ws.Range("F34").Value = "Lavet design af OC6"
$ws.Range("G34").Value = 0.3444444444444445
$ws.Range("H34").Value = 0.43472222222222223

$ws.Range("E35").Value = "Reviewer"
$ws.Range("F35").Value = "Lavet review af OC6 design "
$ws.Range("G35").Value = 0.4548611111111111
$ws.Range("H35").Value = 0.47222222222222227

$ws.Range("F36").Value = "Lavet review af SSD3 "
$ws.Range("G36").Value = 0.47222222222222227
$ws.Range("H36").Value = 0.4770833333333333

$ws.Range("E37").Value = "Implementer"
$ws.Range("F37").Value = "Har Implementeret test case for OC6"
$ws.Range("G37").Value = 0.50347222222222221
$ws.Range("H37").Value = 0.54513888888888895

$ws.Range("F38").Value = "Har Implementeret design for OC6"
$ws.Range("G38").Value = 0.54513888888888895
$ws.Range("H38").Value = 0.57222222222222219

$ws.Range("E39").Value = "Software Architect"
$ws.Range("F39").Value = "Har lavet SSD for UC9 "
$ws.Range("G39").Value = 0.57638888888888895
$ws.Range("H39").Value = 0.58333333333333337

$ws.Range("E40").Value = "Reviewer"
$ws.Range("F40").Value = "Lavet krydstjek af UC6 "
$ws.Range("G40").Value = 0.59375
$ws.Range("H40").Value = 0.62847222222222221

# --- Daily total (row 41) -------------------------------------------------
$ws.Range("I41").Value = 5.0999999999999996

# --- Update the view so the new rows are visible/selected -----------------
$ws.Range("H41").Select()
